$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: was extr1 -> now line7
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

# Row 9: was extr2 -> now line8
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("E9").Value = $true

# Row 10: was extr3 -> now extr1
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

# Row 11: was extr4 -> now extr2
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

# Row 12: was extr5 -> now extr3
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10

# Row 13: was extr6 -> now extr4
$ws.Range("B13").Value = "extr4"
$ws.Range("D13").Value = 8

# Row 14: was extr7 -> now extr5
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11

# Row 15: was extr8 -> now extr6
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11

# New row 16: extr7 (copy A15's format onto A16 so it matches the other
# bold/bordered/centered cells in column A)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

# New row 17: extr8
$ws.Range("A15").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
